$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) — copy formatting (bold / border / centered)
# from the existing last header cell so the added cells reuse the same
# cell style rather than creating a new one, then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record data (Wins/Losses/Ties) for every player row (2-56)
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 103
    $ws.Cells.Item($r, 31).Value = 59
    $ws.Cells.Item($r, 32).Value = 0
}
